# Helper to build an OLE (BGR) color value from R,G,B components,
# mirroring the VBA RGB() function. Must be called as: MyRGB 255 230 153
function MyRGB($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Establish a uniform thin-box border (matching the sheet's existing
#    "border10" style) across the whole new data region A1:N4 by
#    copying the format of an already-bordered cell (B2) onto it.
#    This also seeds the border for the brand new column N and row 4.
# ---------------------------------------------------------------------
$ws.Range("B2").Copy() | Out-Null
$ws.Range("A1:N4").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------------
# 2. Write the new header row.
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "Order Received Data and Time"
$ws.Range("B1").Value = "OrderID"
$ws.Range("C1").Value = "Emp ID-Order Assigned"
$ws.Range("D1").Value = "Assignee_QA"
$ws.Range("E1").Value = "Typist"
$ws.Range("F1").Value = "Typist QC"
$ws.Range("G1").Value = "Client"
$ws.Range("H1").Value = "Lob"
$ws.Range("I1").Value = "Process"
$ws.Range("J1").Value = "Product Name"
$ws.Range("K1").Value = "State"
$ws.Range("L1").Value = "County"
$ws.Range("M1").Value = "Status"
$ws.Range("N1").Value = "Tier"

# ---------------------------------------------------------------------
# 3. Write the new data rows (2-4).
# ---------------------------------------------------------------------
$ws.Range("A2").Value = 45597.0625
$ws.Range("B2").Value = "Acc11-001"
$ws.Range("C2").Value = "SIPL5316"
$ws.Range("D2").Value = "SIPL5688"
$ws.Range("E2").Value = "SIPL5317"
$ws.Range("F2").Value = "SIPL5317"
$ws.Range("G2").Value = "Accurate"
$ws.Range("H2").Value = "Current Owner Search"
$ws.Range("I2").Value = "Search & Typing"
$ws.Range("J2").Value = "One Owner Equity"
$ws.Range("K2").Value = "IA"
$ws.Range("L2").Value = "Adair"
$ws.Range("M2").Value = "WIP"
$ws.Range("N2").Value = ""

$ws.Range("A3").Value = 45613.0625
$ws.Range("B3").Value = "Acc11-002"
$ws.Range("C3").Value = ""
$ws.Range("D3").Value = ""
$ws.Range("E3").Value = "SIPL5317"
$ws.Range("F3").Value = "SIPL5317"
$ws.Range("G3").Value = "Accurate"
$ws.Range("H3").Value = "Equity"
$ws.Range("I3").Value = "Typing"
$ws.Range("J3").Value = "EQ Prop 30yr"
$ws.Range("K3").Value = "IA"
$ws.Range("L3").Value = "Adair"
$ws.Range("M3").Value = "Typing"
$ws.Range("N3").Value = "Typing(T1)"

$ws.Range("A4").Value = 45620.0625
$ws.Range("B4").Value = "Acc11-003"
$ws.Range("C4").Value = ""
$ws.Range("D4").Value = ""
$ws.Range("E4").Value = "SIPL0102"
$ws.Range("F4").Value = "SIPL0103"
$ws.Range("G4").Value = "Accurate"
$ws.Range("H4").Value = "LnV"
$ws.Range("I4").Value = "Typing"
$ws.Range("J4").Value = "Legal and Vesting Report"
$ws.Range("K4").Value = "IA"
$ws.Range("L4").Value = "Appanoose"
$ws.Range("M4").Value = "Typing"
$ws.Range("N4").Value = "Typing(T2)"

# ---------------------------------------------------------------------
# 4. Base font color for the whole region is solid black.
# ---------------------------------------------------------------------
$ws.Range("A1:N4").Font.Color = 0

# ---------------------------------------------------------------------
# 5. Header row formatting: bold black font on a yellow fill with a
#    black pattern color, matching the workbook's header style.
# ---------------------------------------------------------------------
$headerCell = $ws.Range("A1")
$headerCell.Font.Bold = $true
$headerCell.Interior.Pattern = 1
$headerCell.Interior.Color = MyRGB 255 230 153
$headerCell.Interior.PatternColor = MyRGB 0 0 0
$headerCell.Copy() | Out-Null
$ws.Range("B1:N1").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------------
# 6. Date/time column (A2:A4): apply the built-in m/d/yy h:mm format.
# ---------------------------------------------------------------------
$dateCell = $ws.Range("A2")
$dateCell.NumberFormat = "m/d/yy h:mm"
$dateCell.Copy() | Out-Null
$ws.Range("A3:A4").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------------
# 7. County column (L2:L4): smaller font, centered horizontally and
#    vertically.
# ---------------------------------------------------------------------
$countyCell = $ws.Range("L2")
$countyCell.Font.Size = 10
$countyCell.HorizontalAlignment = -4108
$countyCell.VerticalAlignment = -4108
$countyCell.Copy() | Out-Null
$ws.Range("L3:L4").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------------
# 8. F4 gets a border without a left edge (distinguishing it from the
#    rest of the grid).
# ---------------------------------------------------------------------
$ws.Range("F4").Borders.Item(7).LineStyle = 0

# ---------------------------------------------------------------------
# 9. Refresh the used-range dimension and selection to match the
#    post-edit sheet (active cell D7, as in the source workbook).
# ---------------------------------------------------------------------
$ws.Range("D7").Select() | Out-Null
